$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# --- Header / summary cell updates -----------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:48 AM"
$ws.Range("C8").Value = 23183.23
$ws.Range("C9").Value = 105
$ws.Range("C10").Value = "07/21/2025 to 07/27/25"

# --- Insert a new data row for "Point 21" before the Friday TOTAL row ------
# Row 72 currently holds the Friday "TOTAL" row; push it (and everything
# below) down by one row, then populate the freed row 72 with the new line
# item. Row 70 already uses the same zebra-stripe style that row 72 needs,
# so copy its formatting (columns A:I only, to avoid ballooning the used
# range out to XFD) down into the new row.
$ws.Range("A72:I72").EntireRow.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

$ws.Range("A70:I70").Copy()
$ws.Range("A72:I72").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A72").Value = "Point 21"
$ws.Range("B72").Value = "XFR-75-72-120CUT-1B-S"
$ws.Range("C72").Value = "Trans"
$ws.Range("D72").Value = "XFR,75KVA,7.2/12.4kVY,120 CUTOVER,1BG,SS"
$ws.Range("E72").Value = "EA"
$ws.Range("F72").Value = 1
$ws.Range("H72").Value = 203

# --- Update the Friday TOTAL (now on row 73) to include the new line item --
$ws.Range("H73").Value = $ws.Range("H73").Value2 + 203
